$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in order E..T mapped to data array indices 0..15
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{
    2  = @(3, 1, 16.75560066666667, 50.266802, 0.9064438825950116, 0.9064438825950115, 3, 1, 31.96029466666667, 95.88088399999999, 0.3907265741426954, 0.3907265741426953, 535.5139346236631, 4819.625411612968, 0.3541717128989525, 0.3541717128989524)
    3  = @(3, 1, 16.75560066666667, 50.266802, 0.9064438825950116, 0.9064438825950115, 3, 1, 26.34807, 79.04420999999999, 0.3221150253382706, 0.3221150253382706, 441.4777392573799, 3973.29965331642, 0.2919791942098126, 0.2919791942098125)
    4  = @(3, 1, 16.75560066666667, 50.266802, 0.9064438825950116, 0.9064438825950115, 3, 1, 23.48872, 70.46616, 0.287158400519034, 0.287158400519034, 393.5676124911467, 3542.10851242032, 0.2602929754862466, 0.2602929754862465)
    5  = @(3, 1, 0.07688299999999999, 0.230649, 0.004159213770485276, 0.004159213770485276, 3, 1, 31.96029466666667, 95.88088399999999, 0.3907265741426954, 0.3907265741426953, 2.457203334857333, 22.114830013716, 0.001625115347668835, 0.001625115347668835)
    6  = @(3, 1, 0.07688299999999999, 0.230649, 0.004159213770485276, 0.004159213770485276, 3, 1, 26.34807, 79.04420999999999, 0.3221150253382706, 0.3221150253382706, 2.025718665809999, 18.23146799229, 0.001339745249067149, 0.001339745249067149)
    7  = @(3, 1, 0.07688299999999999, 0.230649, 0.004159213770485276, 0.004159213770485276, 3, 1, 23.48872, 70.46616, 0.287158400519034, 0.287158400519034, 1.80588325976, 16.25294933784, 0.001194353173749293, 0.001194353173749293)
    8  = @(3, 1, 1.652500333333333, 4.957501, 0.0893969036345032, 0.08939690363450319, 3, 1, 31.96029466666667, 95.88088399999999, 0.3907265741426954, 0.3907265741426953, 52.81439759009822, 475.329578310884, 0.03492974589607411, 0.0349297458960741)
    9  = @(3, 1, 1.652500333333333, 4.957501, 0.0893969036345032, 0.08939690363450319, 3, 1, 26.34807, 79.04420999999999, 0.3221150253382706, 0.3221150253382706, 43.54019445768999, 391.8617501192099, 0.02879608587939093, 0.02879608587939093)
    10 = @(3, 1, 1.652500333333333, 4.957501, 0.0893969036345032, 0.08939690363450319, 3, 1, 23.48872, 70.46616, 0.287158400519034, 0.287158400519034, 38.81511762957333, 349.33605866616, 0.02567107185903816, 0.02567107185903815)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}
